# Adds two new Mac-Address rows (31 and 32) to the
# master-reg_center_user_machine sheet, mirroring the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31
$ws.Cells.Item(31, 1).Value = 10001
$ws.Cells.Item(31, 2).Value = 110030
$ws.Cells.Item(31, 3).Value = 10030
$ws.Cells.Item(31, 4).Value = "eng"
$ws.Cells.Item(31, 5).Value = $true
$ws.Cells.Item(31, 6).Value = "superadmin"
$ws.Cells.Item(31, 7).Value = "now()"

# Row 32
$ws.Cells.Item(32, 1).Value = 10001
$ws.Cells.Item(32, 2).Value = 110031
$ws.Cells.Item(32, 3).Value = 10031
$ws.Cells.Item(32, 4).Value = "eng"
$ws.Cells.Item(32, 5).Value = $true
$ws.Cells.Item(32, 6).Value = "superadmin"
$ws.Cells.Item(32, 7).Value = "now()"

# Update the view to match what was left selected/scrolled in the author's session.
$ws.Range("C29").Select()
$excel.ActiveWindow.ScrollRow = 25
